# "Generate Report for Archive"
#
# The localization status report was regenerated: rows that used to read
# "Ready for handoff" are now "In Translation", and the Status-ish columns
# that hold that text are narrowed to fit the new (shorter) value.

$wb = $excel.ActiveWorkbook

# --- 1. Update the status text on every sheet that shows it -----------------

$overview = $wb.Worksheets.Item("Overview")
foreach ($col in @("E", "F")) {
    foreach ($row in 2..4) {
        $cell = $overview.Range("$col$row")
        if ($cell.Value2 -eq "Ready for handoff") {
            $cell.Value = "In Translation"
        }
    }
}

foreach ($sheetName in @("zh-cn", "de-de")) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in 2..4) {
        $cell = $ws.Range("C$row")
        if ($cell.Value2 -eq "Ready for handoff") {
            $cell.Value = "In Translation"
        }
    }
}

# --- 2. Narrow the columns that carry that status text ----------------------
# (they were sized for "Ready for handoff" and can shrink for "In Translation")

$overview.Range("E1").ColumnWidth = 12.5
$overview.Range("F1").ColumnWidth = 12.5

$wb.Worksheets.Item("zh-cn").Range("C1").ColumnWidth = 12.5
$wb.Worksheets.Item("de-de").Range("C1").ColumnWidth = 12.5
